$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 7049.5
$ws.Range("J32").Value = 6714.2856
$ws.Range("L32").Value = 6714.2856
$ws.Range("N32").Value = -7366.2856
$ws.Range("H40").Value = 5384.8
$ws.Range("I40").Value = 3641.3333
$ws.Range("K40").Value = 3641.3333
$ws.Range("M40").Value = -3466.3333
$ws.Range("H64").Value = 8583.125
$ws.Range("I64").Value = 4047.1428
$ws.Range("J64").Value = 12111.111
$ws.Range("K64").Value = 4047.1428
$ws.Range("L64").Value = 12111.111
$ws.Range("M64").Value = -3799.1428
$ws.Range("N64").Value = -12607.111
$ws.Range("H67").Value = 8583.125
$ws.Range("I67").Value = 4047.1428
$ws.Range("J67").Value = 12111.111
$ws.Range("K67").Value = 4047.1428
$ws.Range("L67").Value = 12111.111
$ws.Range("M67").Value = -3189.1428
$ws.Range("N67").Value = -13827.111
$ws.Range("H98").Value = 4246.25
$ws.Range("I98").Value = 4229.3022
$ws.Range("J98").Value = 4975
$ws.Range("K98").Value = 4229.3022
$ws.Range("L98").Value = 4975
$ws.Range("M98").Value = -2731.3022
$ws.Range("N98").Value = -7971
$ws.Range("H101").Value = 3039.0715
$ws.Range("J101").Value = 6232.8335
$ws.Range("L101").Value = 18698.5005
$ws.Range("N101").Value = -21942.5005
$ws.Range("H121").Value = 8664.5625
$ws.Range("I121").Value = 7898.7144
$ws.Range("K121").Value = 23696.1432
$ws.Range("M121").Value = -21949.1432
$ws.Range("H122").Value = 4246.25
$ws.Range("I122").Value = 4229.3022
$ws.Range("J122").Value = 4975
$ws.Range("K122").Value = 12687.9066
$ws.Range("L122").Value = 14925
$ws.Range("M122").Value = -10237.9066
$ws.Range("N122").Value = -19825
$ws.Range("H125").Value = 13797.786
$ws.Range("I125").Value = 22651.2
$ws.Range("J125").Value = 8879.223
$ws.Range("K125").Value = 203860.8
$ws.Range("L125").Value = 79913.007
$ws.Range("M125").Value = -201400.8
$ws.Range("N125").Value = -84833.007
$ws.Range("H127").Value = 574.75
$ws.Range("I127").Value = 574.75
$ws.Range("K127").Value = 1724.25
$ws.Range("M127").Value = 3235.75
$ws.Range("H137").Value = 1884.4
$ws.Range("I137").Value = 1307.8889
$ws.Range("J137").Value = 2749.1667
$ws.Range("K137").Value = 3923.6667
$ws.Range("L137").Value = 8247.500100000001
$ws.Range("M137").Value = -1373.6667
$ws.Range("N137").Value = -13347.5001
$ws.Range("H141").Value = 7904
$ws.Range("I141").Value = 7124.375
$ws.Range("K141").Value = 21373.125
$ws.Range("M141").Value = -16193.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 4500
$ws.Range("I63").Value = 4000
$ws.Range("J63").Value = 6000
$ws.Range("K63").Value = 4000
$ws.Range("L63").Value = 6000
$ws.Range("M63").Value = -3314
$ws.Range("N63").Value = -7372
$ws.Range("H66").Value = 4500
$ws.Range("I66").Value = 4000
$ws.Range("J66").Value = 6000
$ws.Range("K66").Value = 20000
$ws.Range("L66").Value = 30000
$ws.Range("M66").Value = -16568
$ws.Range("N66").Value = -36864
$ws.Range("H74").Value = 3253.2888
$ws.Range("J74").Value = 5357.75
$ws.Range("L74").Value = 5357.75
$ws.Range("N74").Value = -7105.75
$ws.Range("H77").Value = 3253.2888
$ws.Range("J77").Value = 5357.75
$ws.Range("L77").Value = 26788.75
$ws.Range("N77").Value = -35524.75
$ws.Range("H88").Value = 1945.1428
$ws.Range("I88").Value = 1567.875
$ws.Range("J88").Value = 2448.1667
$ws.Range("K88").Value = 1567.875
$ws.Range("L88").Value = 2448.1667
$ws.Range("M88").Value = -1161.875
$ws.Range("N88").Value = -3260.1667
$ws.Range("H91").Value = 1945.1428
$ws.Range("I91").Value = 1567.875
$ws.Range("J91").Value = 2448.1667
$ws.Range("K91").Value = 1567.875
$ws.Range("L91").Value = 2448.1667
$ws.Range("M91").Value = -163.875
$ws.Range("N91").Value = -5256.1667
$ws.Range("H102").Value = 3997
$ws.Range("I102").Value = 3139.5715
$ws.Range("K102").Value = 3139.5715
$ws.Range("M102").Value = -1517.5715
$ws.Range("H122").Value = 1288.2307
$ws.Range("I122").Value = 1288.2307
$ws.Range("K122").Value = 3864.6921
$ws.Range("M122").Value = -1414.6921
$ws.Range("H132").Value = 3459.5806
$ws.Range("I132").Value = 2732.5112
$ws.Range("J132").Value = 5384.1763
$ws.Range("K132").Value = 8197.533599999999
$ws.Range("L132").Value = 16152.5289
$ws.Range("M132").Value = -5667.533599999999
$ws.Range("N132").Value = -21212.5289
$ws.Range("H139").Value = 127127.164
$ws.Range("J139").Value = 127127.164
$ws.Range("L139").Value = 127127.164
$ws.Range("N139").Value = -137407.164

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3861.0417
$ws.Range("I134").Value = 3921.2046
$ws.Range("J134").Value = 3199.25
$ws.Range("K134").Value = 11763.6138
$ws.Range("L134").Value = 9597.75
$ws.Range("M134").Value = -9228.613799999999
$ws.Range("N134").Value = -14667.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 661.25
$ws.Range("I2").Value = 381.66666
$ws.Range("K2").Value = 381.66666
$ws.Range("M2").Value = -268.66666
$ws.Range("H4").Value = 25005000
$ws.Range("J4").Value = 10000
$ws.Range("L4").Value = 10000
$ws.Range("N4").Value = -10224
$ws.Range("H10").Value = 549.75
$ws.Range("I10").Value = 549.75
$ws.Range("K10").Value = 549.75
$ws.Range("M10").Value = -410.75
$ws.Range("H38").Value = 4068.5
$ws.Range("I38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("M38").ClearContents()
$ws.Range("H46").Value = 4068.5
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("H58").Value = 3693.3572
$ws.Range("I58").Value = 1784.9474
$ws.Range("K58").Value = 1784.9474
$ws.Range("M58").Value = -1581.9474
$ws.Range("H62").Value = 4852.25
$ws.Range("J62").Value = 6502.5
$ws.Range("L62").Value = 6502.5
$ws.Range("N62").Value = -7750.5
$ws.Range("H65").Value = 4852.25
$ws.Range("J65").Value = 6502.5
$ws.Range("L65").Value = 32512.5
$ws.Range("N65").Value = -38752.5
$ws.Range("H105").Value = 1126.2106
$ws.Range("I105").Value = 911
$ws.Range("K105").Value = 911
$ws.Range("M105").Value = 836
$ws.Range("H132").Value = 2372.3333
$ws.Range("I132").Value = 2133.6365
$ws.Range("K132").Value = 6400.9095
$ws.Range("M132").Value = -3870.9095
$ws.Range("H136").Value = 3693.3572
$ws.Range("I136").Value = 1784.9474
$ws.Range("K136").Value = 5354.8422
$ws.Range("M136").Value = -2804.8422
$ws.Range("H141").Value = 40310.4
$ws.Range("J141").Value = 40310.4
$ws.Range("L141").Value = 40310.4
$ws.Range("N141").Value = -50670.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4014
$ws.Range("J80").Value = 4118.6665
$ws.Range("L80").Value = 4118.6665
$ws.Range("N80").Value = -6114.6665
$ws.Range("H83").Value = 4014
$ws.Range("J83").Value = 4118.6665
$ws.Range("L83").Value = 20593.3325
$ws.Range("N83").Value = -30577.3325
$ws.Range("H122").Value = 2060.3704
$ws.Range("I122").Value = 2114.1765
$ws.Range("J122").Value = 1968.9
$ws.Range("K122").Value = 6342.529500000001
$ws.Range("L122").Value = 5906.700000000001
$ws.Range("M122").Value = -3892.529500000001
$ws.Range("N122").Value = -10806.7

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2567.8333
$ws.Range("J7").Value = 2749.5
$ws.Range("L7").Value = 2749.5
$ws.Range("N7").Value = -2973.5
$ws.Range("H93").Value = 44499.8
$ws.Range("I93").Value = 1250
$ws.Range("K93").Value = 1250
$ws.Range("M93").Value = -2
$ws.Range("H126").Value = 2567.8333
$ws.Range("J126").Value = 2749.5
$ws.Range("L126").Value = 8248.5
$ws.Range("N126").Value = -13188.5
$ws.Range("H132").Value = 6584.778
$ws.Range("I132").Value = 6380.645
$ws.Range("K132").Value = 19141.935
$ws.Range("M132").Value = -16611.935
$ws.Range("H136").Value = 7188.3076
$ws.Range("I136").Value = 7050.4443
$ws.Range("J136").Value = 7498.5
$ws.Range("K136").Value = 21151.3329
$ws.Range("L136").Value = 22495.5
$ws.Range("M136").Value = -18601.3329
$ws.Range("N136").Value = -27595.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 6358.875
$ws.Range("I126").Value = 6420.7
$ws.Range("K126").Value = 19262.1
$ws.Range("M126").Value = -16792.1
